# Weekly update: a new cherry price record was added for
# "Vega Modelo de Temuco" (Early Burlat, Provincia de Limarí) dated
# 2021-11-03 (Excel serial 44503). This pushes all the existing daily
# records (previously rows 11-115) down by one row, and the new record
# takes the freed-up row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11; Excel shifts rows 11..115 down to
# 12..116 automatically (carries the date-format style on column D along).
$ws.Rows(11).Insert()

# Populate the newly inserted row 11 with the new record.
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Vega Modelo de Temuco"
$ws.Range("C11").Value = "La Araucanía"
$ws.Range("D11").Value = 44503
$ws.Range("E11").Value = 9
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100103
$ws.Range("H11").Value = "Frutos de hueso (carozo)"
$ws.Range("I11").Value = 100103001
$ws.Range("J11").Value = "Cereza"
$ws.Range("K11").Value = "Early Burlat"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 20000
$ws.Range("Q11").Value = "`$/bandeja 10 kilos"
$ws.Range("R11").Value = "Provincia de Limarí"
$ws.Range("S11").Value = 2000
$ws.Range("T11").Value = 10
